# "gran crest start book"
# Adds a new 2013 "Grancrest Start Book 1" (Fujimi Shobo) row to the
# checklist, right after the existing 2013 rulebook entries, and fixes
# several English-title typos ("Grand Rest"/"Grand Crest" -> "Grancrest"
# / "Grancest") in the process.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new row ------------------------------------------------
# Existing row 4 (2014, rulebook 2) and everything below shifts down one.
$ws.Rows.Item(4).Insert()

$ws.Range("A4").Value = 2013
$ws.Range("B4").Value = "グランクレストスタートブック 1"
$ws.Range("C4").Value = "Granrest Start Book 1"
$ws.Range("D4").Value = "Fujimi Shobo"
$ws.Range("E4").Value = "start_book.jpg"
$ws.Range("F4").Value = "rulebook"

# --- Fix existing English titles (column C) -----------------------------
$ws.Range("C2").Value  = "Grancrest RPG Rulebook 1"
$ws.Range("C5").Value  = "Grancest RPG Rulebook 2"
$ws.Range("C6").Value  = "Grancrest RPG Data Book Advanced Force"
$ws.Range("C14").Value = "Grancrest RPG Supplement Advanced Rulebook"
$ws.Range("C15").Value = "Record of Grancrest War Data Book"

# --- Refresh the sort state / dimension so it covers the new row 20 ----
$sortObj = $ws.Sort
$sortObj.SortFields.Clear() | Out-Null
$sortObj.SortFields.Add($ws.Range("A2:A20")) | Out-Null
$sortObj.SetRange($ws.Range("A2:F20")) | Out-Null
$sortObj.Header = 0
$sortObj.Apply() | Out-Null

# --- Match the final saved selection ------------------------------------
$ws.Range("D5").Select() | Out-Null
